$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap data between rows 236/239 and 237/238 ---
# (On 2024-04-01 the 4 already-played matches from 2023-10-06 (id 234-237) were
#  re-synced from source and came back in a different order; id/date/odds etc.
#  for rows 236 and 239 swap places, and likewise for rows 237 and 238.)
# Row 236
$ws.Range("B236").Formula = "6836277"
$ws.Range("C236").Value = "Romania Liga I"
$ws.Range("D236").Value = "Romania Liga I"
$ws.Range("E236").Formula = "45359.625"
$ws.Range("F236").Value = "CFR Cluj"
$ws.Range("G236").Value = "AFC Hermannstadt"
$ws.Range("H236").Formula = "1"
$ws.Range("I236").Formula = "0"
$ws.Range("J236").Value = "H"
$ws.Range("K236").Formula = "1.7"
$ws.Range("L236").Formula = "3.4"
$ws.Range("M236").Formula = "5"
$ws.Range("N236").Formula = "1.65"
$ws.Range("O236").Formula = "3.5"
$ws.Range("P236").Formula = "5.25"
$ws.Range("Q236").Formula = "-0.75"
$ws.Range("R236").Formula = "1.85"
$ws.Range("S236").Formula = "2"
$ws.Range("T236").Formula = "2.25"
$ws.Range("U236").Formula = "1.875"
$ws.Range("V236").Formula = "1.975"
$ws.Range("W236").Formula = "0.6499999999999999"
$ws.Range("X236").Formula = "-1"
$ws.Range("Y236").Formula = "-1"
$ws.Range("Z236").Formula = "0.425"
$ws.Range("AA236").Formula = "-0.5"
$ws.Range("AB236").Formula = "-1"
$ws.Range("AC236").Formula = "0.9750000000000001"

# Row 237
$ws.Range("B237").Formula = "6870268"
$ws.Range("C237").Value = "Romania Liga I"
$ws.Range("D237").Value = "Romania Liga I"
$ws.Range("E237").Formula = "45359.625"
$ws.Range("F237").Value = "Petrolul Ploiesti"
$ws.Range("G237").Value = "ACS Sepsi"
$ws.Range("H237").Formula = "1"
$ws.Range("I237").Formula = "2"
$ws.Range("J237").Value = "A"
$ws.Range("K237").Formula = "2.8"
$ws.Range("L237").Formula = "3"
$ws.Range("M237").Formula = "2.55"
$ws.Range("N237").Formula = "3"
$ws.Range("O237").Formula = "3.2"
$ws.Range("P237").Formula = "2.3"
$ws.Range("Q237").Formula = "0.25"
$ws.Range("R237").Formula = "1.85"
$ws.Range("S237").Formula = "2"
$ws.Range("T237").Formula = "2.25"
$ws.Range("U237").Formula = "1.875"
$ws.Range("V237").Formula = "1.975"
$ws.Range("W237").Formula = "-1"
$ws.Range("X237").Formula = "-1"
$ws.Range("Y237").Formula = "1.3"
$ws.Range("Z237").Formula = "-1"
$ws.Range("AA237").Formula = "1"
$ws.Range("AB237").Formula = "0.875"
$ws.Range("AC237").Formula = "-1"

# Row 238
$ws.Range("B238").Formula = "6865915"
$ws.Range("C238").Value = "Romania Liga I"
$ws.Range("D238").Value = "Romania Liga I"
$ws.Range("E238").Formula = "45359.625"
$ws.Range("F238").Value = "FC Voluntari"
$ws.Range("G238").Value = "Universitatea Cluj"
$ws.Range("H238").Formula = "0"
$ws.Range("I238").Formula = "0"
$ws.Range("J238").Value = "D"
$ws.Range("K238").Formula = "3.5"
$ws.Range("L238").Formula = "3.25"
$ws.Range("M238").Formula = "2.05"
$ws.Range("N238").Formula = "3.4"
$ws.Range("O238").Formula = "3.1"
$ws.Range("P238").Formula = "2.15"
$ws.Range("Q238").Formula = "0.25"
$ws.Range("R238").Formula = "1.975"
$ws.Range("S238").Formula = "1.875"
$ws.Range("T238").Formula = "2.25"
$ws.Range("U238").Formula = "2.05"
$ws.Range("V238").Formula = "1.75"
$ws.Range("W238").Formula = "-1"
$ws.Range("X238").Formula = "2.1"
$ws.Range("Y238").Formula = "-1"
$ws.Range("Z238").Formula = "0.4875"
$ws.Range("AA238").Formula = "-0.5"
$ws.Range("AB238").Formula = "-1"
$ws.Range("AC238").Formula = "0.75"

# Row 239
$ws.Range("B239").Formula = "6861095"
$ws.Range("C239").Value = "Romania Liga I"
$ws.Range("D239").Value = "Romania Liga I"
$ws.Range("E239").Formula = "45359.625"
$ws.Range("F239").Value = "FC Botosani"
$ws.Range("G239").Value = "Farul Constanta"
$ws.Range("H239").Formula = "0"
$ws.Range("I239").Formula = "0"
$ws.Range("J239").Value = "D"
$ws.Range("K239").Formula = "3.75"
$ws.Range("L239").Formula = "3.4"
$ws.Range("M239").Formula = "1.909"
$ws.Range("N239").Formula = "3.1"
$ws.Range("O239").Formula = "3"
$ws.Range("P239").Formula = "2.375"
$ws.Range("Q239").Formula = "0.25"
$ws.Range("R239").Formula = "1.775"
$ws.Range("S239").Formula = "2.1"
$ws.Range("T239").Formula = "2"
$ws.Range("U239").Formula = "1.8"
$ws.Range("V239").Formula = "2.05"
$ws.Range("W239").Formula = "-1"
$ws.Range("X239").Formula = "2"
$ws.Range("Y239").Formula = "-1"
$ws.Range("Z239").Formula = "0.3875"
$ws.Range("AA239").Formula = "-0.5"
$ws.Range("AB239").Formula = "-1"
$ws.Range("AC239").Formula = "1.05"

# --- Row 256: result data refreshed / corrected (same fixture, new odds) ---
$ws.Range("A256").Formula = "254"
$ws.Range("B256").Formula = "7951779"
$ws.Range("C256").Value = "Romania Liga I"
$ws.Range("D256").Value = "Romania Liga I"
$ws.Range("E256").Formula = "45382.33333333334"
$ws.Range("F256").Value = "FC U Craiova 1948"
$ws.Range("G256").Value = "Otelul Galati"
$ws.Range("H256").Formula = "1"
$ws.Range("I256").Formula = "2"
$ws.Range("J256").Value = "A"
$ws.Range("K256").Formula = "2.3"
$ws.Range("L256").Formula = "3.2"
$ws.Range("M256").Formula = "3.2"
$ws.Range("N256").Formula = "2.15"
$ws.Range("O256").Formula = "3.25"
$ws.Range("P256").Formula = "3.4"
$ws.Range("Q256").Formula = "-0.25"
$ws.Range("R256").Formula = "1.875"
$ws.Range("S256").Formula = "1.975"
$ws.Range("T256").Formula = "2.25"
$ws.Range("U256").Formula = "2.05"
$ws.Range("V256").Formula = "1.8"
$ws.Range("W256").Formula = "-1"
$ws.Range("X256").Formula = "-1"
$ws.Range("Y256").Formula = "2.4"
$ws.Range("Z256").Formula = "-1"
$ws.Range("AA256").Formula = "0.9750000000000001"
$ws.Range("AB256").Formula = "1.05"
$ws.Range("AC256").Formula = "-1"

# --- New rows 257-263: newly added fixtures ---
# Row 257
$ws.Cells.Item(255, 1).Copy($ws.Cells.Item(257, 1))
$ws.Cells.Item(255, 5).Copy($ws.Cells.Item(257, 5))
$ws.Range("A257").Formula = "255"
$ws.Range("B257").Formula = "7951748"
$ws.Range("C257").Value = "Romania Liga I"
$ws.Range("D257").Value = "Romania Liga I"
$ws.Range("E257").Formula = "45382.625"
$ws.Range("F257").Value = "Farul Constanta"
$ws.Range("G257").Value = "FCSB"
$ws.Range("H257").Formula = "0"
$ws.Range("I257").Formula = "1"
$ws.Range("J257").Value = "A"
$ws.Range("K257").Formula = "3.6"
$ws.Range("L257").Formula = "3.3"
$ws.Range("M257").Formula = "2"
$ws.Range("N257").Formula = "3.6"
$ws.Range("O257").Formula = "3.4"
$ws.Range("P257").Formula = "2"
$ws.Range("Q257").Formula = "0.5"
$ws.Range("R257").Formula = "1.8"
$ws.Range("S257").Formula = "2.05"
$ws.Range("T257").Formula = "2.5"
$ws.Range("U257").Formula = "2"
$ws.Range("V257").Formula = "1.85"
$ws.Range("W257").Formula = "-1"
$ws.Range("X257").Formula = "-1"
$ws.Range("Y257").Formula = "1"
$ws.Range("Z257").Formula = "-1"
$ws.Range("AA257").Formula = "1.05"
$ws.Range("AB257").Formula = "-1"
$ws.Range("AC257").Formula = "0.8500000000000001"

# Row 258
$ws.Cells.Item(255, 1).Copy($ws.Cells.Item(258, 1))
$ws.Cells.Item(255, 5).Copy($ws.Cells.Item(258, 5))
$ws.Range("A258").Formula = "256"
$ws.Range("B258").Formula = "7951780"
$ws.Range("C258").Value = "Romania Liga I"
$ws.Range("D258").Value = "Romania Liga I"
$ws.Range("E258").Formula = "45383.60416666666"
$ws.Range("F258").Value = "Dinamo Bucharest"
$ws.Range("G258").Value = "Petrolul Ploiesti"
$ws.Range("H258").Formula = "1"
$ws.Range("I258").Formula = "1"
$ws.Range("J258").Value = "D"
$ws.Range("K258").Formula = "2.3"
$ws.Range("L258").Formula = "3"
$ws.Range("M258").Formula = "3.4"
$ws.Range("N258").Formula = "2.4"
$ws.Range("O258").Formula = "2.9"
$ws.Range("P258").Formula = "3.3"
$ws.Range("Q258").Formula = "-0.25"
$ws.Range("R258").Formula = "2"
$ws.Range("S258").Formula = "1.85"
$ws.Range("T258").Formula = "2"
$ws.Range("U258").Formula = "2.05"
$ws.Range("V258").Formula = "1.8"
$ws.Range("W258").Formula = "-1"
$ws.Range("X258").Formula = "1.9"
$ws.Range("Y258").Formula = "-1"
$ws.Range("Z258").Formula = "-0.5"
$ws.Range("AA258").Formula = "0.425"
$ws.Range("AB258").Formula = "0"
$ws.Range("AC258").Formula = "-0"

# Row 259
$ws.Cells.Item(255, 1).Copy($ws.Cells.Item(259, 1))
$ws.Cells.Item(255, 5).Copy($ws.Cells.Item(259, 5))
$ws.Range("A259").Formula = "257"
$ws.Range("B259").Formula = "7951783"
$ws.Range("C259").Value = "Romania Liga I"
$ws.Range("D259").Value = "Romania Liga I"
$ws.Range("E259").Formula = "45387.47916666666"
$ws.Range("F259").Value = "ACS UTA Batrana Doamna"
$ws.Range("G259").Value = "FC Botosani"
$ws.Range("K259").Formula = "1.95"
$ws.Range("L259").Formula = "3.3"
$ws.Range("M259").Formula = "4"
$ws.Range("N259").Formula = "1.95"
$ws.Range("O259").Formula = "3.3"
$ws.Range("P259").Formula = "4"
$ws.Range("Q259").Formula = "-0.5"
$ws.Range("R259").Formula = "2"
$ws.Range("S259").Formula = "1.85"
$ws.Range("T259").Formula = "2.25"
$ws.Range("U259").Formula = "1.925"
$ws.Range("V259").Formula = "1.925"
$ws.Range("W259").Formula = "0"
$ws.Range("X259").Formula = "0"
$ws.Range("Y259").Formula = "0"
$ws.Range("Z259").Formula = "0"
$ws.Range("AA259").Formula = "0"

# Row 260
$ws.Cells.Item(255, 1).Copy($ws.Cells.Item(260, 1))
$ws.Cells.Item(255, 5).Copy($ws.Cells.Item(260, 5))
$ws.Range("A260").Formula = "258"
$ws.Range("B260").Formula = "7951753"
$ws.Range("C260").Value = "Romania Liga I"
$ws.Range("D260").Value = "Romania Liga I"
$ws.Range("E260").Formula = "45387.60416666666"
$ws.Range("F260").Value = "ACS Sepsi"
$ws.Range("G260").Value = "Farul Constanta"
$ws.Range("K260").Formula = "2.15"
$ws.Range("L260").Formula = "3.25"
$ws.Range("M260").Formula = "3.4"
$ws.Range("N260").Formula = "2.15"
$ws.Range("O260").Formula = "3.25"
$ws.Range("P260").Formula = "3.4"
$ws.Range("Q260").Formula = "-0.25"
$ws.Range("R260").Formula = "1.875"
$ws.Range("S260").Formula = "1.975"
$ws.Range("T260").Formula = "2.25"
$ws.Range("U260").Formula = "1.825"
$ws.Range("V260").Formula = "2.025"
$ws.Range("W260").Formula = "0"
$ws.Range("X260").Formula = "0"
$ws.Range("Y260").Formula = "0"
$ws.Range("Z260").Formula = "0"
$ws.Range("AA260").Formula = "0"

# Row 261
$ws.Cells.Item(255, 1).Copy($ws.Cells.Item(261, 1))
$ws.Cells.Item(255, 5).Copy($ws.Cells.Item(261, 5))
$ws.Range("A261").Formula = "259"
$ws.Range("B261").Formula = "7951752"
$ws.Range("C261").Value = "Romania Liga I"
$ws.Range("D261").Value = "Romania Liga I"
$ws.Range("E261").Formula = "45388.60416666666"
$ws.Range("F261").Value = "Rapid Bucuresti"
$ws.Range("G261").Value = "CFR Cluj"
$ws.Range("K261").Formula = "2.8"
$ws.Range("L261").Formula = "3.1"
$ws.Range("M261").Formula = "2.6"
$ws.Range("N261").Formula = "2.8"
$ws.Range("O261").Formula = "3.1"
$ws.Range("P261").Formula = "2.6"
$ws.Range("Q261").Formula = "0"
$ws.Range("R261").Formula = "2"
$ws.Range("S261").Formula = "1.85"
$ws.Range("T261").Formula = "2.25"
$ws.Range("U261").Formula = "1.875"
$ws.Range("V261").Formula = "1.975"
$ws.Range("W261").Formula = "0"
$ws.Range("X261").Formula = "0"
$ws.Range("Y261").Formula = "0"
$ws.Range("Z261").Formula = "0"
$ws.Range("AA261").Formula = "0"

# Row 262
$ws.Cells.Item(255, 1).Copy($ws.Cells.Item(262, 1))
$ws.Cells.Item(255, 5).Copy($ws.Cells.Item(262, 5))
$ws.Range("A262").Formula = "260"
$ws.Range("B262").Formula = "7951784"
$ws.Range("C262").Value = "Romania Liga I"
$ws.Range("D262").Value = "Romania Liga I"
$ws.Range("E262").Formula = "45389.35416666666"
$ws.Range("F262").Value = "AFC Hermannstadt"
$ws.Range("G262").Value = "Universitatea Cluj"
$ws.Range("K262").Formula = "2.4"
$ws.Range("L262").Formula = "2.875"
$ws.Range("M262").Formula = "3.3"
$ws.Range("N262").Formula = "2.4"
$ws.Range("O262").Formula = "2.875"
$ws.Range("P262").Formula = "3.3"
$ws.Range("Q262").Formula = "-0.25"
$ws.Range("R262").Formula = "2.05"
$ws.Range("S262").Formula = "1.8"
$ws.Range("T262").Formula = "2"
$ws.Range("U262").Formula = "2.025"
$ws.Range("V262").Formula = "1.825"
$ws.Range("W262").Formula = "0"
$ws.Range("X262").Formula = "0"
$ws.Range("Y262").Formula = "0"
$ws.Range("Z262").Formula = "0"
$ws.Range("AA262").Formula = "0"

# Row 263
$ws.Cells.Item(255, 1).Copy($ws.Cells.Item(263, 1))
$ws.Cells.Item(255, 5).Copy($ws.Cells.Item(263, 5))
$ws.Range("A263").Formula = "261"
$ws.Range("B263").Formula = "7951751"
$ws.Range("C263").Value = "Romania Liga I"
$ws.Range("D263").Value = "Romania Liga I"
$ws.Range("E263").Formula = "45389.60416666666"
$ws.Range("F263").Value = "FCSB"
$ws.Range("G263").Value = "CS U Craiova"
$ws.Range("K263").Formula = "1.85"
$ws.Range("L263").Formula = "3.4"
$ws.Range("M263").Formula = "4.2"
$ws.Range("N263").Formula = "1.85"
$ws.Range("O263").Formula = "3.4"
$ws.Range("P263").Formula = "4.2"
$ws.Range("Q263").Formula = "-0.5"
$ws.Range("R263").Formula = "1.875"
$ws.Range("S263").Formula = "1.975"
$ws.Range("T263").Formula = "2.5"
$ws.Range("U263").Formula = "1.925"
$ws.Range("V263").Formula = "1.925"
$ws.Range("W263").Formula = "0"
$ws.Range("X263").Formula = "0"
$ws.Range("Y263").Formula = "0"
$ws.Range("Z263").Formula = "0"
$ws.Range("AA263").Formula = "0"
